# Curation fix for the kelp frond counts workbook: the "year" column (A) on
# every site sheet was stored as a full four-digit year ("2016", "2017",
# "2018"); it should instead hold the curated 2-digit numeric year code
# (17, 18, 19). Applies to rows 2-4 (the data rows below the header) on
# every worksheet, skipping rows that don't exist on a given sheet (e.g.
# sheets with only two data rows).

$wb = $excel.ActiveWorkbook

$yearMap = @{
    "2016" = 17
    "2017" = 18
    "2018" = 19
}

foreach ($ws in $wb.Worksheets) {
    foreach ($row in 2..4) {
        $cell = $ws.Cells.Item($row, 1)
        $current = $cell.Value()
        if ($current -eq $null) {
            continue
        }
        $key = [string]$current
        if ($yearMap.ContainsKey($key)) {
            $cell.Value = $yearMap[$key]
        }
    }
}
